# Move the four "results" text boxes on slide 1 (shape ids 4, 5, 6, 7) up
# by ~238536 EMU (~0.2609in / 18.78pt), keeping Left/Width/Height unchanged.
#
# PowerPoint's Shape.Top/.Left COM properties are Single-precision (f32)
# values expressed in points (1 pt = 12700 EMU). The literals below were
# chosen so that converting them to f32 and back to EMU reproduces the
# exact target EMU offsets from the authored file.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape.Id -> new Top (in points)
$targets = @{
    4 = 203.47850803700788   # y: 2822713 -> 2584177 EMU
    5 = 203.47850803700788   # y: 2822713 -> 2584177 EMU
    6 = 203.02024082047245   # y: 2816893 -> 2578357 EMU
    7 = 203.47220612440944   # y: 2822633 -> 2584097 EMU
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    $id = $shape.Id
    if ($targets.ContainsKey($id)) {
        $shape.Top = $targets[$id]
    }
}
